# Applies the ELM-1NA -> ELM-2NA schedule swap for rows 18-21 (columns B:F)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Sandro-Acionamentos Elétricos, Pedro Bispo-Manut. Elétrica, Sandro-Lógica de Programação, Cleidson-Sistemas digitais]"
$ws.Range("C18").Value = "Euclides-Gestão Integr"
$ws.Range("D18").Value = "[Leonardo-Manut. Mecânica, Rogério-Processos de Usinagem 2, Weslei-Metrologia, Anderson-Processos de Usinagem 1]"
$ws.Range("F18").Value = "[Rachel-Tecnologia dos Materiais., Rachel-Tecnologia dos Materiais.]"

$ws.Range("B19").Value = "[Sandro-Acionamentos Elétricos, Pedro Bispo-Manut. Elétrica, Sandro-Lógica de Programação, Cleidson-Sistemas digitais]"
$ws.Range("C19").Value = "Guilherme-Máquinas Térmicas e de Fl"
$ws.Range("D19").Value = "[Leonardo-Manut. Mecânica, Rogério-Processos de Usinagem 2, Weslei-Metrologia, Anderson-Processos de Usinagem 1]"
$ws.Range("E19").Value = "Gilberto-Tecnologias Mecâni"

$ws.Range("B20").Value = "[Sandro-Acionamentos Elétricos, Pedro Bispo-Manut. Elétrica, Sandro-Lógica de Programação, Cleidson-Sistemas digitais]"
$ws.Range("C20").Value = "Allan Cupertino-Circuitos Elétrico"
$ws.Range("D20").Value = "[Leonardo-Manut. Mecânica, Rogério-Processos de Usinagem 2, Weslei-Metrologia, Anderson-Processos de Usinagem 1]"
$ws.Range("E20").Value = "Guilherme-Máquinas Térmicas e de Fl"
$ws.Range("F20").Value = "[Suzanny-Des. Bas. Mec., Suzanny-Des. Bas. Mec.]"

$ws.Range("B21").Value = "[Sandro-Acionamentos Elétricos, Pedro Bispo-Manut. Elétrica, Sandro-Lógica de Programação, Cleidson-Sistemas digitais]"
$ws.Range("C21").Value = "Allan Cupertino-Circuitos Elétrico"
$ws.Range("D21").Value = "[Leonardo-Manut. Mecânica, Rogério-Processos de Usinagem 2, Weslei-Metrologia, Anderson-Processos de Usinagem 1]"
$ws.Range("F21").Value = "[Suzanny-Des. Bas. Mec., Suzanny-Des. Bas. Mec.]"
